$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3762.8333
$ws.Range("I69").Value = 3790.6
$ws.Range("K69").Value = 11371.8
$ws.Range("M69").Value = -10497.8

$ws.Range("H72").Value = 3762.8333
$ws.Range("I72").Value = 3790.6
$ws.Range("K72").Value = 34115.4
$ws.Range("M72").Value = -29747.4

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").ClearContents()

$ws.Range("H115").Value = 1033.75
$ws.Range("I115").Value = 378.33334
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 1135.00002
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = 431.9999800000001
$ws.Range("N115").Value = -12134

$ws.Range("H139").Value = 44978
$ws.Range("J139").Value = 45743.848
$ws.Range("L139").Value = 45743.848
$ws.Range("N139").Value = -56023.848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2685.4
$ws.Range("I61").Value = 2141.3462
$ws.Range("J61").Value = 4257.1113
$ws.Range("K61").Value = 2141.3462
$ws.Range("L61").Value = 4257.1113
$ws.Range("M61").Value = -1929.3462
$ws.Range("N61").Value = -4681.1113

$ws.Range("H88").Value = 3003.5
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 3003.5
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H97").Value = 920.86957
$ws.Range("I97").Value = 844.3125
$ws.Range("J97").Value = 1095.8572
$ws.Range("K97").Value = 844.3125
$ws.Range("L97").Value = 1095.8572
$ws.Range("M97").Value = -348.3125
$ws.Range("N97").Value = -2087.8572

$ws.Range("H136").Value = 2685.4
$ws.Range("I136").Value = 2141.3462
$ws.Range("J136").Value = 4257.1113
$ws.Range("K136").Value = 6424.0386
$ws.Range("L136").Value = 12771.3339
$ws.Range("M136").Value = -3874.0386
$ws.Range("N136").Value = -17871.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 32856.19
$ws.Range("J126").Value = 32856.19
$ws.Range("L126").Value = 32856.19
$ws.Range("N126").Value = -42736.19

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6788.7676
$ws.Range("I31").Value = 1528.4783
$ws.Range("J31").Value = 12838.1
$ws.Range("K31").Value = 1528.4783
$ws.Range("L31").Value = 12838.1
$ws.Range("M31").Value = -1233.4783
$ws.Range("N31").Value = -13428.1

$ws.Range("H34").Value = 6788.7676
$ws.Range("I34").Value = 1528.4783
$ws.Range("J34").Value = 12838.1
$ws.Range("K34").Value = 1528.4783
$ws.Range("L34").Value = 12838.1
$ws.Range("M34").Value = -1326.4783
$ws.Range("N34").Value = -13242.1

$ws.Range("H122").Value = 1945.0555
$ws.Range("I122").Value = 1670.3334
$ws.Range("K122").Value = 5011.0002
$ws.Range("M122").Value = -2561.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2291.8572
$ws.Range("I69").Value = 452.2
$ws.Range("J69").Value = 2866.75
$ws.Range("K69").Value = 1356.6
$ws.Range("L69").Value = 8600.25
$ws.Range("M69").Value = -545.5999999999999
$ws.Range("N69").Value = -10222.25

$ws.Range("H72").Value = 2291.8572
$ws.Range("I72").Value = 452.2
$ws.Range("J72").Value = 2866.75
$ws.Range("K72").Value = 4069.8
$ws.Range("L72").Value = 25800.75
$ws.Range("M72").Value = -13.79999999999973
$ws.Range("N72").Value = -33912.75

$ws.Range("H74").Value = 3950
$ws.Range("J74").Value = 3950
$ws.Range("L74").Value = 11850
$ws.Range("N74").Value = -13972

$ws.Range("H77").Value = 3950
$ws.Range("J77").Value = 3950
$ws.Range("L77").Value = 35550
$ws.Range("N77").Value = -46158

$ws.Range("H80").Value = 9823.333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9823.333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29469.999
$ws.Range("N80").Value = -31341.999
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 9823.333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9823.333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 88409.997
$ws.Range("N83").Value = -97769.997
$ws.Range("M83").ClearContents()

$ws.Range("H113").Value = 892.5854
$ws.Range("I113").Value = 643.3889
$ws.Range("J113").Value = 1087.6086
$ws.Range("K113").Value = 1930.1667
$ws.Range("L113").Value = 3262.8258
$ws.Range("M113").Value = 239.8332999999998
$ws.Range("N113").Value = -7602.825800000001

$ws.Range("H122").Value = 3685.6667
$ws.Range("I122").Value = 386.0625
$ws.Range("J122").Value = 6791.1763
$ws.Range("K122").Value = 3474.5625
$ws.Range("L122").Value = 61120.5867
$ws.Range("M122").Value = -1024.5625
$ws.Range("N122").Value = -66020.5867

$ws.Range("H129").Value = 1588.4584
$ws.Range("I129").Value = 391.42856
$ws.Range("K129").Value = 1174.28568
$ws.Range("M129").Value = 3825.71432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 85260.75
$ws.Range("I113").Value = 251370
$ws.Range("J113").Value = 2206.125
$ws.Range("K113").Value = 251370
$ws.Range("L113").Value = 2206.125
$ws.Range("M113").Value = -249200
$ws.Range("N113").Value = -6546.125

$ws.Range("H139").Value = 66842
$ws.Range("J139").Value = 66842
$ws.Range("L139").Value = 66842
$ws.Range("N139").Value = -77122

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4433.3335
$ws.Range("I61").Value = 4675
$ws.Range("J61").Value = 4240
$ws.Range("K61").Value = 4675
$ws.Range("L61").Value = 4240
$ws.Range("M61").Value = -4473
$ws.Range("N61").Value = -4644

$ws.Range("H113").Value = 4433.3335
$ws.Range("I113").Value = 4675
$ws.Range("J113").Value = 4240
$ws.Range("K113").Value = 4675
$ws.Range("L113").Value = 4240
$ws.Range("M113").Value = -2505
$ws.Range("N113").Value = -8580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 41539.5
$ws.Range("J46").Value = 41539.5
$ws.Range("L46").Value = 41539.5
$ws.Range("N46").Value = -42001.5

$ws.Range("H81").Value = 4360.636
$ws.Range("I81").Value = 5928.1665
$ws.Range("J81").Value = 2479.6
$ws.Range("K81").Value = 11856.333
$ws.Range("L81").Value = 4959.2
$ws.Range("M81").Value = -10795.333
$ws.Range("N81").Value = -7081.2

$ws.Range("H84").Value = 4360.636
$ws.Range("I84").Value = 5928.1665
$ws.Range("J84").Value = 2479.6
$ws.Range("K84").Value = 59281.665
$ws.Range("L84").Value = 24796
$ws.Range("M84").Value = -53977.665
$ws.Range("N84").Value = -35404

$ws.Range("H107").Value = 833.5
$ws.Range("J107").Value = 499.66666
$ws.Range("L107").Value = 1498.99998
$ws.Range("N107").Value = -5338.999980000001

$ws.Range("H134").Value = 41539.5
$ws.Range("J134").Value = 41539.5
$ws.Range("L134").Value = 124618.5
$ws.Range("N134").Value = -129688.5

$ws.Range("H138").Value = 57543.2
$ws.Range("J138").Value = 61929
$ws.Range("L138").Value = 61929
$ws.Range("N138").Value = -72209
